$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '64.730.71'
$ws.Range("E2").Value = '  +2.30%  '
$ws.Range("D3").Value = '3.372.41'
$ws.Range("E3").Value = '  +2.08%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '557.99'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +2.32%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '174.91'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.74%  '
$ws.Range("E7").Value = '  +2.13%  '
$ws.Range("D8").Value = '3.362.12'
$ws.Range("E8").Value = '  +2.06%  '
$ws.Range("E9").Value = '  +0.00%  '
$ws.Range("E10").Value = '  +8.90%  '
$ws.Range("E11").Value = '  +3.49%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '53.93'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.67%  '
$ws.Range("E13").Value = '  +4.27%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.08'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.50%  '
$ws.Range("D15").Value = '3.903.94'
$ws.Range("E15").Value = '  +2.15%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '18.24'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +1.24%  '
$ws.Range("E17").Value = '  +2.39%  '
$ws.Range("D18").Value = '3.382.32'
$ws.Range("E18").Value = '  +2.43%  '
$ws.Range("D19").Value = '64.520.26'
$ws.Range("E19").Value = '  +2.14%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.76'
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '0.991'
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '459.87'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +6.57%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '4.86'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +9.03%  '
$ws.Range("E24").Value = '  +2.02%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '86.22'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.96%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '13.57'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +2.89%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '2.94'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +8.67%  '
$ws.Range("E28").Value = '  +2.06%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.72'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +1.51%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '30.61'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +5.27%  '
$ws.Range("E31").Value = '  +5.45%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.43'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +1.25%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '571.06'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.00%  '
$ws.Range("E34").Value = '  +5.35%  '
$ws.Range("E35").Value = '  +1.95%  '
$ws.Range("E36").Value = '  +0.03%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.62'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +6.58%  '
$ws.Range("E38").Value = '  -4.15%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '35.38'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +1.36%  '
$ws.Range("D40").Value = '0.0₃0742'
$ws.Range("E40").Value = '  +1.09%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.369'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +1.88%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.999'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.04%  '
$ws.Range("D43").Value = '3.073.25'
$ws.Range("E43").Value = '  -1.12%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.81'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.38%  '
$ws.Range("E45").Value = '  +3.65%  '
$ws.Range("E46").Value = '  +4.67%  '
$ws.Range("E47").Value = '  +1.35%  '
$ws.Range("E48").Value = '  -3.21%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.58'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +0.26%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '138.19'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +4.34%  '
$ws.Range("E51").Value = '  +2.55%  '
